$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "34÷7=4, 6"
$t.Cell(1, 2).Range.Text = "32÷9=3, 5"
$t.Cell(1, 3).Range.Text = "96÷9=10, 6"
$t.Cell(1, 4).Range.Text = "10÷5=2, 0"
$t.Cell(1, 5).Range.Text = "19÷8=2, 3"
$t.Cell(5, 1).Range.Text = "18÷8=2, 2"
$t.Cell(5, 2).Range.Text = "33÷3=11, 0"
$t.Cell(5, 3).Range.Text = "54÷9=6, 0"
$t.Cell(5, 4).Range.Text = "49÷7=7, 0"
$t.Cell(5, 5).Range.Text = "77÷5=15, 2"
$t.Cell(9, 1).Range.Text = "86÷9=9, 5"
$t.Cell(9, 2).Range.Text = "69÷2=34, 1"
$t.Cell(9, 3).Range.Text = "66÷9=7, 3"
$t.Cell(9, 4).Range.Text = "74÷4=18, 2"
$t.Cell(9, 5).Range.Text = "95÷4=23, 3"
$t.Cell(13, 1).Range.Text = "81÷3=27, 0"
$t.Cell(13, 2).Range.Text = "10÷4=2, 2"
$t.Cell(13, 3).Range.Text = "24÷9=2, 6"
$t.Cell(13, 4).Range.Text = "59÷4=14, 3"
$t.Cell(13, 5).Range.Text = "36÷4=9, 0"
$t.Cell(17, 1).Range.Text = "77÷9=8, 5"
$t.Cell(17, 2).Range.Text = "46÷8=5, 6"
$t.Cell(17, 3).Range.Text = "30÷7=4, 2"
$t.Cell(17, 4).Range.Text = "87÷4=21, 3"
$t.Cell(17, 5).Range.Text = "56÷2=28, 0"
